$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write a literal text value (e.g. "1.72 %") into a cell without
# Excel's automatic percentage-number coercion kicking in. We park the text
# behind a temporary "=""...""" text formula, then copy/paste-special
# (values only) the cell onto itself, which collapses the formula down to a
# plain shared-string cell exactly like the other literal "NN %" cells
# already in this sheet (no leftover formula, no leftover cell style).
function Set-LiteralText($addr, $text) {
    $rng = $ws.Range($addr)
    $escaped = $text -replace '"', '""'
    $rng.Formula = '="' + $escaped + '"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}

# --- "Contenu du stage" block (rows 16-23) ---------------------------------
# Only C#, COBOL, ASSEMBLEUR and ANDROID actually gain students; the rest of
# the block (C++, JEE, DELPHI, PHP5) stays at 0 / "0 %" so those rows are
# left untouched.

# Row 16 : C# -> 1 etudiant / 1.72 %
$ws.Range("E16").Value = 1
Set-LiteralText "G16" "1.72 %"

# Row 17 : COBOL -> 52 etudiants / 89.66 %
$ws.Range("E17").Value = 52
Set-LiteralText "G17" "89.66 %"

# Row 19 : ASSEMBLEUR -> 2 etudiants / 3.45 %
$ws.Range("E19").Value = 2
Set-LiteralText "G19" "3.45 %"

# Row 20 : ANDROID -> 3 etudiants / 5.17 %
$ws.Range("E20").Value = 3
Set-LiteralText "G20" "5.17 %"
